$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update dSF (column F) values to reflect repulled data / recalculated mean
$ws.Range("F2").Value = 6
$ws.Range("F3").Value = -3
$ws.Range("F4").Value = -5
$ws.Range("F6").Value = -5
$ws.Range("F8").Value = 3
